$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 = "I0" and J1 = "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the style of the existing header cells (e.g. H1: bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Fill I and J columns for data rows 2-37
# I column is always 1, J column mirrors the H column value
for ($r = 2; $r -le 37; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
